# Updated symbol list on Wed Feb 15 15:37:36 UTC 2023 with GitHub Actions
#
# This script updates the "Price" (column D) and "Volume(1h)" (column E)
# values in the crypto tracker worksheet to the latest scraped figures.
# Values are written as literal text (matching the workbook's existing
# inline-string storage) by forcing a "Text" number format on each cell
# before assigning it, so Excel does not reinterpret numeric- or
# percent-looking strings as numbers.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet


$ws.Range("D2").NumberFormat = "@"
$ws.Range("D2").Value = "302.77"

$ws.Range("E2").NumberFormat = "@"
$ws.Range("E2").Value = "1.94%"

$ws.Range("D3").NumberFormat = "@"
$ws.Range("D3").Value = "43.30"

$ws.Range("E3").NumberFormat = "@"
$ws.Range("E3").Value = "6.67%"

$ws.Range("D4").NumberFormat = "@"
$ws.Range("D4").Value = "5.080"

$ws.Range("E4").NumberFormat = "@"
$ws.Range("E4").Value = "0.50%"

$ws.Range("D5").NumberFormat = "@"
$ws.Range("D5").Value = "0.07696"

$ws.Range("E5").NumberFormat = "@"
$ws.Range("E5").Value = "3.42%"

$ws.Range("D6").NumberFormat = "@"
$ws.Range("D6").Value = "1.623"

$ws.Range("E6").NumberFormat = "@"
$ws.Range("E6").Value = "3.11%"

$ws.Range("D7").NumberFormat = "@"
$ws.Range("D7").Value = "1.045"

$ws.Range("E7").NumberFormat = "@"
$ws.Range("E7").Value = "11.71%"

$ws.Range("D8").NumberFormat = "@"
$ws.Range("D8").Value = "0.1253"

$ws.Range("E8").NumberFormat = "@"
$ws.Range("E8").Value = "3.60%"

$ws.Range("D9").NumberFormat = "@"
$ws.Range("D9").Value = "0.1859"

$ws.Range("E9").NumberFormat = "@"
$ws.Range("E9").Value = "2.81%"

$ws.Range("D10").NumberFormat = "@"
$ws.Range("D10").Value = "0.09183"

$ws.Range("E10").NumberFormat = "@"
$ws.Range("E10").Value = "4.13%"

$ws.Range("D11").NumberFormat = "@"
$ws.Range("D11").Value = "0.04167"

$ws.Range("E11").NumberFormat = "@"
$ws.Range("E11").Value = "-2.31%"

$ws.Range("D12").NumberFormat = "@"
$ws.Range("D12").Value = "0.1048"

$ws.Range("E12").NumberFormat = "@"
$ws.Range("E12").Value = "-0.09%"

$ws.Range("D13").NumberFormat = "@"
$ws.Range("D13").Value = "0.001271"

$ws.Range("E13").NumberFormat = "@"
$ws.Range("E13").Value = "-1.68%"

$ws.Range("D14").NumberFormat = "@"
$ws.Range("D14").Value = "0.005750"

$ws.Range("E14").NumberFormat = "@"
$ws.Range("E14").Value = "-2.99%"

$ws.Range("E15").NumberFormat = "@"
$ws.Range("E15").Value = "1,901.46%"

$ws.Range("D16").NumberFormat = "@"
$ws.Range("D16").Value = "3.340"

$ws.Range("E16").NumberFormat = "@"
$ws.Range("E16").Value = "-0.51%"

$ws.Range("D17").NumberFormat = "@"
$ws.Range("D17").Value = "4.412"

$ws.Range("E17").NumberFormat = "@"
$ws.Range("E17").Value = "1.29%"

$ws.Range("E18").NumberFormat = "@"
$ws.Range("E18").Value = "-1.98%"

$ws.Range("D19").NumberFormat = "@"
$ws.Range("D19").Value = "0.3354"

$ws.Range("E19").NumberFormat = "@"
$ws.Range("E19").Value = "1.42%"

$ws.Range("D20").NumberFormat = "@"
$ws.Range("D20").Value = "8.650"

$ws.Range("E20").NumberFormat = "@"
$ws.Range("E20").Value = "9.55%"

$ws.Range("E21").NumberFormat = "@"
$ws.Range("E21").Value = "-0.70%"

$ws.Range("E22").NumberFormat = "@"
$ws.Range("E22").Value = "8.31%"

$ws.Range("D23").NumberFormat = "@"
$ws.Range("D23").Value = "0.04169"

$ws.Range("E23").NumberFormat = "@"
$ws.Range("E23").Value = "3.85%"

$ws.Range("D24").NumberFormat = "@"
$ws.Range("D24").Value = "0.001284"

$ws.Range("E24").NumberFormat = "@"
$ws.Range("E24").Value = "1.54%"

$ws.Range("D25").NumberFormat = "@"
$ws.Range("D25").Value = "0.004473"

$ws.Range("E25").NumberFormat = "@"
$ws.Range("E25").Value = "15.46%"

$ws.Range("D26").NumberFormat = "@"
$ws.Range("D26").Value = "0.0001347"

$ws.Range("E26").NumberFormat = "@"
$ws.Range("E26").Value = "9.81%"

$ws.Range("D38").NumberFormat = "@"
$ws.Range("D38").Value = "0.02459"

$ws.Range("E38").NumberFormat = "@"
$ws.Range("E38").Value = "2.63%"

$ws.Range("D39").NumberFormat = "@"
$ws.Range("D39").Value = "0.05280"

$ws.Range("E39").NumberFormat = "@"
$ws.Range("E39").Value = "1.45%"

$ws.Range("D40").NumberFormat = "@"
$ws.Range("D40").Value = "0.005967"

$ws.Range("E40").NumberFormat = "@"
$ws.Range("E40").Value = "0.22%"

$ws.Range("D41").NumberFormat = "@"
$ws.Range("D41").Value = "0.007683"

$ws.Range("E41").NumberFormat = "@"
$ws.Range("E41").Value = "-1.06%"

$ws.Range("D42").NumberFormat = "@"
$ws.Range("D42").Value = "0.1345"

$ws.Range("E42").NumberFormat = "@"
$ws.Range("E42").Value = "1.88%"

$ws.Range("D43").NumberFormat = "@"
$ws.Range("D43").Value = "0.007374"

$ws.Range("E43").NumberFormat = "@"
$ws.Range("E43").Value = "0.40%"

$ws.Range("D44").NumberFormat = "@"
$ws.Range("D44").Value = "0.007558"

$ws.Range("E44").NumberFormat = "@"
$ws.Range("E44").Value = "5.14%"

$ws.Range("E45").NumberFormat = "@"
$ws.Range("E45").Value = "2.03%"

$ws.Range("D46").NumberFormat = "@"
$ws.Range("D46").Value = "0.00006713"

$ws.Range("E46").NumberFormat = "@"
$ws.Range("E46").Value = "7.17%"

$ws.Range("D47").NumberFormat = "@"
$ws.Range("D47").Value = "0.00000000748"

$ws.Range("E47").NumberFormat = "@"
$ws.Range("E47").Value = "0.03%"

$ws.Range("D48").NumberFormat = "@"
$ws.Range("D48").Value = "0.04447"

$ws.Range("E48").NumberFormat = "@"
$ws.Range("E48").Value = "-2.31%"

$ws.Range("E49").NumberFormat = "@"
$ws.Range("E49").Value = "0.27%"

$ws.Range("D50").NumberFormat = "@"
$ws.Range("D50").Value = "0.00002096"

$ws.Range("E50").NumberFormat = "@"
$ws.Range("E50").Value = "0.03%"

$ws.Range("D51").NumberFormat = "@"
$ws.Range("D51").Value = "0.0001996"

$ws.Range("E51").NumberFormat = "@"
$ws.Range("E51").Value = "0.03%"
